# Apply cryptos list update (price + 1h volume/change columns)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.903.80'
$ws.Range("E2").Value = '  -0.56%  '
$ws.Range("D3").Value = '2.814.42'
$ws.Range("E3").Value = '  +0.89%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '353.05'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.41%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '111.40'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.94%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.564'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.71%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E9").Value = '  +2.82%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.55'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -5.89%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0854'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.25%  '
$ws.Range("E12").Value = '  +0.26%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.80'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.56%  '
$ws.Range("E14").Value = '  +0.28%  '
$ws.Range("D15").Value = '3.252.92'
$ws.Range("E15").Value = '  +0.79%  '
$ws.Range("D16").Value = '2.807.28'
$ws.Range("E16").Value = '  +0.73%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.925'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +4.20%  '
$ws.Range("D18").Value = '51.698.00'
$ws.Range("E18").Value = '  -0.71%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.54'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +6.55%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.11'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.24%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.33'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.18%  '
$ws.Range("D22").Value = '0.0₃0990'
$ws.Range("E22").Value = '  +0.94%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.22'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.14%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '267.86'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.88%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.80'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.37%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.83'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.69%  '
$ws.Range("E27").Value = '  +0.02%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.26'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.07%  '
$ws.Range("E29").Value = '  +0.60%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0492'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +21.11%  '
$ws.Range("E31").Value = '  -0.41%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '52.57'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.67%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '34.38'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.20%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.90'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.11%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.53'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +11.34%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0845'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.46%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.999'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.10%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.24'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.32%  '
$ws.Range("E39").Value = '  -4.16%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.27'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.18%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.117'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.71%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '127.08'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.50%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '22.91'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.50%  '
$ws.Range("E44").Value = '  -2.22%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.47'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -8.64%  '
$ws.Range("D46").Value = '2.083.45'
$ws.Range("E46").Value = '  +0.51%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.33'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.49%  '
$ws.Range("E49").Value = '  +6.78%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.977'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +7.60%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '9.05'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.03%  '
